$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# New configuration rows describing the System1 reports download/upload workflow
$ws.Range("A9").Value = "System1_ReportsPath"
$ws.Range("B9").Value = "reports/"
$ws.Range("C9").Value = "Url path to get the reports page page on System1 website"

$ws.Range("A10").Value = "System1_ReportsDownload"
$ws.Range("B10").Value = "download"
$ws.Range("C10").Value = "Url path to get the reports download page on System1 website"

$ws.Range("A11").Value = "System1_DownloadFilePath"
$ws.Range("B11").Value = "C:\Users\dsembiante\Downloads"
$ws.Range("C11").Value = "Folder reports are downloaded into"

$ws.Range("A12").Value = "System1_MoveFilePath"
$ws.Range("B12").Value = "C:\Users\dsembiante\OneDrive - Deloitte (O365D)\UiPathAcademyArcitect\Reports"
$ws.Range("C12").Value = "Folder reports are moved into"

$ws.Range("A13").Value = "System1_ReportNameFilter"
$ws.Range("B13").Value = "Report-"
$ws.Range("C13").Value = "Part of the file path for downlaoded reports"

$ws.Range("A14").Value = "System1_ReportsUpload"
$ws.Range("B14").Value = "upload"
$ws.Range("C14").Value = "Url path to get the reports upload page on System1 website"

# Match the wrap-text style used on other description cells in column C
$ws.Range("C9").WrapText = $true
$ws.Range("C11").WrapText = $true
$ws.Range("C13").WrapText = $true

# Update the active selection on the sheet
$ws.Range("C16").Select() | Out-Null
